$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.841.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7755"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3135"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07344"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7656"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.458"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.871.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.195"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.854.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007848"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.163"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.127.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1573"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.424"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.038"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.453"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.543"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.476"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05568"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.079"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7569"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9956"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01931"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.146.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4447"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8515"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.899"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.125"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.826"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.496"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
